$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.749.55'
$ws.Range("E2").Value = '  -1.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.683.70'
$ws.Range("E3").Value = '  -1.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '221.35'
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5230'
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.009'
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06557'
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2591'
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.29'
$ws.Range("E10").Value = '  -2.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07711'
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.714.98'
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.924.51'
$ws.Range("E13").Value = '  -0.90%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.369'
$ws.Range("E14").Value = '  -4.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5647'
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8068'
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.82'
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.876.10'
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.29'
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.530'
$ws.Range("E21").Value = '  -2.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.20'
$ws.Range("E22").Value = '  -2.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.938'
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.009'
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.95'
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.734'
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1184'
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.085'
$ws.Range("E28").Value = '  -2.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.94'
$ws.Range("E29").Value = '  -2.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05279'
$ws.Range("E30").Value = '  -1.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.278'
$ws.Range("E31").Value = '  -0.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.388'
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.271'
$ws.Range("E33").Value = '  -4.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.606'
$ws.Range("E34").Value = '  -1.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.788'
$ws.Range("E35").Value = '  -3.07%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.393'
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9344'
$ws.Range("E37").Value = '  -1.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5801'
$ws.Range("E38").Value = '  -0.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.171.78'
$ws.Range("E39").Value = '  +12.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01621'
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.008'
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.718'
$ws.Range("E42").Value = '  -2.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8327'
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.95'
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.832.14'
$ws.Range("E45").Value = '  -0.92%  '
$ws.Range("E46").Value = '  -4.04%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4508'
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '56.22'
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.007'
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.026'
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05188'
$ws.Range("E51").Value = '  -0.97%  '
